$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Loading Details Name" column (P) with a sample row value,
# matching the style of the existing header row / data row.

# Header cell P7 - copy formatting from an existing simple header cell (A7)
# then set its text.
$headerSrc = $ws.Range("A7")
$headerDst = $ws.Range("P7")
$headerSrc.Copy()
$headerDst.PasteSpecial(-4122)  # xlPasteFormats
$headerDst.Value = "Loading Details Name"

# Data cell P8 - copy formatting from the adjacent "loop devices" data cell
# (K8, which uses the grey/bordered data style) then set its text.
$dataSrc = $ws.Range("K8")
$dataDst = $ws.Range("P8")
$dataSrc.Copy()
$dataDst.PasteSpecial(-4122)  # xlPasteFormats
$dataDst.Value = "Main Processor 24V (A)"

$ws.Application.CutCopyMode = $false

# Resize the new column to fit its content.
$ws.Columns.Item(16).AutoFit()

# Update page setup to portrait orientation.
$ws.PageSetup.Orientation = 1

# Move the view / selection to the newly added column, mirroring the
# author scrolling right to review the new data.
$win = $excel.ActiveWindow
$win.ScrollColumn = 13
$win.ScrollRow = 1
$ws.Columns.Item(16).Select() | Out-Null
